$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells to Text format first so Excel does not
# auto-convert numeric-looking strings (e.g. "2.340", "0.09080") into
# actual numbers and silently drop formatting (trailing zeros, etc.).
$priceCells = @("D2","D3","D4","D5","D6","D7","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "24.641.36"
$ws.Range("E2").Value = "  +3.61%  "
$ws.Range("D3").Value = "1.699.70"
$ws.Range("E3").Value = "  +2.46%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "317.56"
$ws.Range("E5").Value = "  +3.31%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "0.3948"
$ws.Range("E7").Value = "  +1.64%  "
$ws.Range("E8").Value = "  +2.63%  "
$ws.Range("E9").Value = "  +9.15%  "
$ws.Range("D10").Value = "55.11"
$ws.Range("E10").Value = "  +13.85%  "
$ws.Range("D11").Value = "1.001"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "0.08813"
$ws.Range("E12").Value = "  +2.60%  "
$ws.Range("D13").Value = "7.278"
$ws.Range("E13").Value = "  +8.94%  "
$ws.Range("D14").Value = "23.44"
$ws.Range("E14").Value = "  +3.62%  "
$ws.Range("D15").Value = "0.00001333"
$ws.Range("E15").Value = "  +2.39%  "
$ws.Range("D16").Value = "7.644"
$ws.Range("E16").Value = "  +6.51%  "
$ws.Range("D17").Value = "1.702.06"
$ws.Range("E17").Value = "  +2.59%  "
$ws.Range("D18").Value = "101.29"
$ws.Range("E18").Value = "  +1.68%  "
$ws.Range("D19").Value = "0.07089"
$ws.Range("E19").Value = "  +4.85%  "
$ws.Range("D20").Value = "19.76"
$ws.Range("E20").Value = "  +4.55%  "
$ws.Range("D21").Value = "6.933"
$ws.Range("E21").Value = "  +4.69%  "
$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "14.18"
$ws.Range("E23").Value = "  +3.40%  "
$ws.Range("D24").Value = "24.641.33"
$ws.Range("E24").Value = "  +3.67%  "
$ws.Range("D25").Value = "2.981"
$ws.Range("E25").Value = "  +10.45%  "
$ws.Range("D26").Value = "2.340"
$ws.Range("E26").Value = "  +1.22%  "
$ws.Range("D27").Value = "22.44"
$ws.Range("D28").Value = "159.72"
$ws.Range("E28").Value = "  +1.83%  "
$ws.Range("D29").Value = "5.262"
$ws.Range("E29").Value = "  +1.76%  "
$ws.Range("D30").Value = "134.02"
$ws.Range("E30").Value = "  +3.64%  "
$ws.Range("D31").Value = "7.612"
$ws.Range("E31").Value = "  +19.54%  "
$ws.Range("B32").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C32").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D32").Value = "1.885.94"
$ws.Range("E32").Value = "  +2.33%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "1.105"
$ws.Range("E33").Value = "  -1.61%  "
$ws.Range("D34").Value = "7.398"
$ws.Range("E34").Value = "  +14.33%  "
$ws.Range("D35").Value = "0.08577"
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("D36").Value = "11.31"
$ws.Range("E36").Value = "  +10.49%  "
$ws.Range("D37").Value = "0.2769"
$ws.Range("E37").Value = "  +5.74%  "
$ws.Range("D38").Value = "1.956"
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("D39").Value = "14.77"
$ws.Range("E39").Value = "  +3.06%  "
$ws.Range("D40").Value = "0.02823"
$ws.Range("E40").Value = "  +12.38%  "
$ws.Range("D41").Value = "0.09080"
$ws.Range("E41").Value = "  +4.00%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.7762"
$ws.Range("E42").Value = "  +3.67%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "1.465"
$ws.Range("E43").Value = "  +1.74%  "
$ws.Range("D44").Value = "0.7288"
$ws.Range("E44").Value = "  +4.48%  "
$ws.Range("D45").Value = "15.65"
$ws.Range("E45").Value = "  +6.11%  "
$ws.Range("D46").Value = "2.520"
$ws.Range("E46").Value = "  +6.73%  "
$ws.Range("D47").Value = "4.221"
$ws.Range("E47").Value = "  +4.25%  "
$ws.Range("D48").Value = "1.389"
$ws.Range("E48").Value = "  +21.45%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").Value = "141.86"
$ws.Range("E50").Value = "  +2.04%  "
$ws.Range("D51").Value = "0.08042"
$ws.Range("E51").Value = "  +3.67%  "
